$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new value looks numeric but must stay plain text
#     (matches source data formatting, e.g. trailing zeros / fixed decimals) ---
$textCells = @(
    @{ Ref = 'D5'; Value = '332.59' }
    @{ Ref = 'D7'; Value = '0.4710' }
    @{ Ref = 'D8'; Value = '0.3956' }
    @{ Ref = 'D9'; Value = '47.88' }
    @{ Ref = 'D10'; Value = '0.08048' }
    @{ Ref = 'D11'; Value = '1.032' }
    @{ Ref = 'D12'; Value = '22.17' }
    @{ Ref = 'D14'; Value = '5.979' }
    @{ Ref = 'D15'; Value = '7.121' }
    @{ Ref = 'D16'; Value = '1.006' }
    @{ Ref = 'D17'; Value = '87.10' }
    @{ Ref = 'D18'; Value = '0.00001046' }
    @{ Ref = 'D19'; Value = '0.06669' }
    @{ Ref = 'D20'; Value = '17.23' }
    @{ Ref = 'D21'; Value = '1.003' }
    @{ Ref = 'D23'; Value = '5.529' }
    @{ Ref = 'D25'; Value = '2.305' }
    @{ Ref = 'D27'; Value = '159.44' }
    @{ Ref = 'D29'; Value = '2.103' }
    @{ Ref = 'D30'; Value = '5.579' }
    @{ Ref = 'D31'; Value = '121.74' }
    @{ Ref = 'D32'; Value = '0.9822' }
    @{ Ref = 'D33'; Value = '0.09525' }
    @{ Ref = 'D34'; Value = '1.447' }
    @{ Ref = 'D35'; Value = '3.598' }
    @{ Ref = 'D36'; Value = '5.356' }
    @{ Ref = 'D37'; Value = '0.06123' }
    @{ Ref = 'D38'; Value = '0.02262' }
    @{ Ref = 'D39'; Value = '1.229' }
    @{ Ref = 'D40'; Value = '8.146' }
    @{ Ref = 'D41'; Value = '0.6019' }
    @{ Ref = 'D42'; Value = '1.003' }
    @{ Ref = 'D43'; Value = '0.1904' }
    @{ Ref = 'D44'; Value = '10.30' }
    @{ Ref = 'D45'; Value = '1.268' }
    @{ Ref = 'D46'; Value = '0.5704' }
    @{ Ref = 'D47'; Value = '12.25' }
    @{ Ref = 'D48'; Value = '1.950' }
    @{ Ref = 'D49'; Value = '3.398' }
    @{ Ref = 'D50'; Value = '0.06919' }
    @{ Ref = 'D51'; Value = '114.06' }
)
foreach ($tc in $textCells) {
    $rng = $ws.Range($tc.Ref)
    $rng.NumberFormat = "@"
    $rng.Value = $tc.Value
    $rng.Style = "Normal"
}

# --- Remaining cells (coin names, links, already-non-numeric prices, % deltas) ---
$ws.Range('D2').Value = '27.733.55'
$ws.Range('E2').Value = '  +1.63%  '
$ws.Range('D3').Value = '1.877.42'
$ws.Range('E3').Value = '  +1.25%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('E5').Value = '  +2.66%  '
$ws.Range('E6').Value = '  +0.23%  '
$ws.Range('E7').Value = '  +3.53%  '
$ws.Range('E8').Value = '  +2.17%  '
$ws.Range('E9').Value = '  -0.81%  '
$ws.Range('E10').Value = '  +1.37%  '
$ws.Range('E11').Value = '  +1.67%  '
$ws.Range('E12').Value = '  +3.71%  '
$ws.Range('D13').Value = '1.886.50'
$ws.Range('E13').Value = '  +1.75%  '
$ws.Range('E14').Value = '  +1.36%  '
$ws.Range('E15').Value = '  -0.28%  '
$ws.Range('E16').Value = '  +0.48%  '
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('E17').Value = '  +1.41%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('E18').Value = '  +1.81%  '
$ws.Range('E19').Value = '  +1.43%  '
$ws.Range('E20').Value = '  +0.96%  '
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('D22').Value = '27.752.38'
$ws.Range('E22').Value = '  +1.67%  '
$ws.Range('E23').Value = '  +0.50%  '
$ws.Range('E24').Value = '  +1.04%  '
$ws.Range('E25').Value = '  +0.76%  '
$ws.Range('D26').Value = '2.103.50'
$ws.Range('E26').Value = '  +1.36%  '
$ws.Range('E27').Value = '  +3.67%  '
$ws.Range('E29').Value = '  +2.13%  '
$ws.Range('E30').Value = '  +2.01%  '
$ws.Range('E31').Value = '  +0.43%  '
$ws.Range('E32').Value = '  +4.88%  '
$ws.Range('E33').Value = '  +2.06%  '
$ws.Range('E34').Value = '  -0.46%  '
$ws.Range('E35').Value = '  +0.36%  '
$ws.Range('E36').Value = '  +1.61%  '
$ws.Range('E37').Value = '  +1.98%  '
$ws.Range('E39').Value = '  +1.03%  '
$ws.Range('E40').Value = '  +1.08%  '
$ws.Range('E41').Value = '  +1.79%  '
$ws.Range('B42').Value = 'Frax'
$ws.Range('C42').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('E42').Value = '  +0.24%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('E43').Value = '  +0.86%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('E44').Value = '  +1.23%  '
$ws.Range('B45').Value = 'WEMIXTOKEN'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('E45').Value = '  -1.17%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('E46').Value = '  +1.83%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E47').Value = '  +2.29%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('E48').Value = '  +1.78%  '
$ws.Range('B49').Value = 'PancakeSwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('E49').Value = '  +0.77%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('E50').Value = '  +2.69%  '
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('E51').Value = '  +4.91%  '
